$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$eValues = @{
    2 = "Alias Changed for Network"
    3 = "Alias Changed for Network"
    4 = "Network Added to Base Service"
    5 = "Network Removed from Database"
    6 = "Network Removed from Database"
    7 = "New Network Added to Database"
    8 = "Network Removed from Database"
    9 = "New Network Added to Database"
    10 = "New Network Added to Database"
    11 = "New Network Added to Database"
    12 = "New Network Added to Database"
    13 = "New Network Added to Database"
    14 = "New Network Added to Database"
    15 = "New Network Added to Database"
    16 = "Network Removed from Database"
    17 = "Network Removed from Database"
    18 = "New Network Added to Database"
    19 = "Network Removed from Database"
    20 = "Network Added to Base Service"
    21 = "Network Added to Add-On Package"
    22 = "Network Added to Add-On Package"
    23 = "Network Added to Add-On Package"
    24 = "Name of Add-On Package Changed"
    25 = "Name of Add-On Package Changed"
    26 = "Name of Add-On Package Changed"
    27 = "Network Added to Base Service"
    28 = "Network Removed from Base Service"
    29 = "Network Removed from Base Service"
    30 = "Network Added to Base Service"
    31 = "Network Removed from Base Service"
    32 = "Network Removed from Base Service"
    33 = "Network Removed from Base Service"
    34 = "Network Added to Base Service"
    35 = "Network Removed from Base Service"
    36 = "Network Added to Base Service"
    37 = "Network Added to Base Service"
    38 = "Network Added to Base Service"
    39 = "Network Added to Base Service"
    40 = "Network Removed from Base Service"
    41 = "Network Added to Base Service"
    42 = "Network Removed from Base Service"
    43 = "Network Removed from Base Service"
    44 = "Network Added to Base Service"
    45 = "Network Added to Base Service"
    46 = "Network Added to Base Service"
    47 = "Network Added to Base Service"
    48 = "Network Added to Base Service"
    49 = "Network Added to Base Service"
    50 = "Network Added to Base Service"
    51 = "Network Added to Base Service"
    52 = "Network Added to Base Service"
    53 = "Network Added to Base Service"
    54 = "Network Added to Base Service"
    55 = "Network Added to Base Service"
    56 = "Network Added to Base Service"
    57 = "Network Added to Base Service"
    58 = "Network Removed from Base Service"
    59 = "Network Removed from Base Service"
    60 = "Network Added to Base Service"
    61 = "Network Removed from Base Service"
    62 = "Network Removed from Base Service"
    63 = "Network Removed from Base Service"
    64 = "Network Added to Base Service"
    65 = "Network Added to Base Service"
    66 = "Network Removed from Base Service"
    67 = "Network Removed from Base Service"
    68 = "Network Removed from Base Service"
    69 = "Network Removed from Base Service"
    70 = "Network Added to Base Service"
    71 = "Network Added to Base Service"
    72 = "Network Added to Base Service"
    73 = "Network Added to Base Service"
    74 = "Network Removed from Base Service"
    75 = "Network Removed from Base Service"
    76 = "Network Added to Base Service"
    77 = "Network Added to Base Service"
    78 = "Network Removed from Base Service"
    79 = "Network Removed from Base Service"
    80 = "Network Removed from Base Service"
    81 = "Network Added to Base Service"
    82 = "Network Added to Base Service"
    83 = "Network Added to Base Service"
    84 = "Network Added to Base Service"
    85 = "Network Added to Base Service"
    86 = "Network Removed from Base Service"
    87 = "Network Added to Base Service"
    88 = "Network Added to Base Service"
    89 = "Network Added to Base Service"
    90 = "Network Removed from Base Service"
    91 = "Network Added to Base Service"
    92 = "Network Added to Base Service"
    93 = "Network Added to Base Service"
    94 = "Network Added to Base Service"
    95 = "Network Added to Base Service"
    96 = "Network Added to Base Service"
    97 = "Network Added to Base Service"
    98 = "Network Added to Base Service"
    99 = "Network Added to Base Service"
    100 = "Network Added to Base Service"
    101 = "Network Added to Base Service"
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $eValues[$row]
}

# Clear and re-apply the autofilter to refresh its uid and drop the stale
# embedded sortState from the previous sort/filter pass.
$ws.Range("A1:E101").AutoFilter() | Out-Null
$ws.Range("A1:E101").AutoFilter() | Out-Null

# User reselected cell B1 last.
$ws.Range("B1").Select() | Out-Null
